$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted as row 13 (most recent
# date first), pushing the existing rows 13-44 down to 14-45 and growing the
# sheet's used range from A1:R44 to A1:R45.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new Cilantro price record for
# Terminal Hortofrutícola Agro Chillán (Ñuble).
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 44687
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 100112040
$ws.Cells.Item(13, 7).Value = "Cilantro"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 300
$ws.Cells.Item(13, 11).Value = 550
$ws.Cells.Item(13, 12).Value = 600
$ws.Cells.Item(13, 13).Value = 575
$ws.Cells.Item(13, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(13, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(13, 16).Value = 575
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(13, 4).NumberFormat = $ws.Cells.Item(14, 4).NumberFormat
